$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.922.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.288.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.630'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.87%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.648'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.38%  '
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.631.21'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.868'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.285.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.811.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.57%  '
$ws.Range("E25").Value = '  -2.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.21%  '
$ws.Range("E34").Value = '  -3.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0819'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.11'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.126'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.73%  '
$ws.Range("E40").Value = '  -4.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.85%  '
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("E43").Value = '  -5.01%  '
$ws.Range("E44").Value = '  +7.40%  '
$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.63%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.79%  '
$ws.Range("E47").Value = '  -2.49%  '
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.67%  '
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("E51").Value = '  -2.03%  '
